# Auto-generated edit script for cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.970.36"
$ws.Range("E2").Value = "  +4.78%  "
$ws.Range("D3").Value = "2.332.11"
$ws.Range("E3").Value = "  +2.62%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'521.58"
$ws.Range("E5").Value = "  +4.69%  "
$ws.Range("D6").Value = "'135.25"
$ws.Range("E6").Value = "  +4.47%  "
$ws.Range("D7").Value = "'0.993"
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("D8").Value = "'0.539"
$ws.Range("E8").Value = "  +2.48%  "
$ws.Range("D9").Value = "2.363.97"
$ws.Range("E9").Value = "  +3.47%  "
$ws.Range("E10").Value = "  +8.77%  "
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("D12").Value = "'5.20"
$ws.Range("E12").Value = "  +6.27%  "
$ws.Range("E13").Value = "  +2.79%  "
$ws.Range("E14").Value = "  +4.39%  "
$ws.Range("D15").Value = "2.753.57"
$ws.Range("E15").Value = "  +2.86%  "
$ws.Range("D16").Value = "57.036.65"
$ws.Range("E16").Value = "  +4.94%  "
$ws.Range("E17").Value = "  +4.92%  "
$ws.Range("D18").Value = "2.342.02"
$ws.Range("E18").Value = "  +2.83%  "
$ws.Range("D19").Value = "'10.63"
$ws.Range("E19").Value = "  +3.32%  "
$ws.Range("D20").Value = "'4.30"
$ws.Range("E20").Value = "  +3.51%  "
$ws.Range("D21").Value = "'323.71"
$ws.Range("E21").Value = "  +6.20%  "
$ws.Range("D22").Value = "'6.67"
$ws.Range("E22").Value = "  +5.08%  "
$ws.Range("D23").Value = "'0.996"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("D25").Value = "'0.989"
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("E26").Value = "  +6.91%  "
$ws.Range("D27").Value = "'7.80"
$ws.Range("E27").Value = "  +5.91%  "
$ws.Range("D28").Value = "'172.37"
$ws.Range("E28").Value = "  -1.39%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "0.0₃0750"
$ws.Range("E29").Value = "  +5.68%  "
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "'1.22"
$ws.Range("E30").Value = "  +12.45%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").Value = "'6.35"
$ws.Range("E31").Value = "  +5.58%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.71"
$ws.Range("E32").Value = "  +5.69%  "
$ws.Range("D33").Value = "'18.47"
$ws.Range("E33").Value = "  +3.34%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").Value = "'0.958"
$ws.Range("E35").Value = "  +0.79%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'1.27"
$ws.Range("E36").Value = "  +5.17%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "'0.990"
$ws.Range("E37").Value = "  -0.69%  "
$ws.Range("D38").Value = "'4.06"
$ws.Range("E38").Value = "  +8.74%  "
$ws.Range("D39").Value = "'1.53"
$ws.Range("E39").Value = "  +8.70%  "
$ws.Range("E40").Value = "  +4.30%  "
$ws.Range("D41").Value = "'0.384"
$ws.Range("E41").Value = "  +1.99%  "
$ws.Range("D42").Value = "'139.95"
$ws.Range("E42").Value = "  +11.64%  "
$ws.Range("D43").Value = "'3.63"
$ws.Range("E43").Value = "  +6.91%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").Value = "'5.23"
$ws.Range("E44").Value = "  +2.88%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "'279.32"
$ws.Range("E45").Value = "  +13.98%  "
$ws.Range("D46").Value = "'0.0512"
$ws.Range("E46").Value = "  +3.83%  "
$ws.Range("E47").Value = "  +3.77%  "
$ws.Range("D48").Value = "'0.567"
$ws.Range("E48").Value = "  +3.28%  "
$ws.Range("D49").Value = "'0.0217"
$ws.Range("E49").Value = "  +5.72%  "
$ws.Range("D50").Value = "'0.383"
$ws.Range("E50").Value = "  +1.94%  "
$ws.Range("D51").Value = "'17.15"
$ws.Range("E51").Value = "  +5.38%  "
